$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused "Uberball" activity name and insert four new blank
# rows so the existing text rows (orig rows 2-5) move down to rows 3,5,7,9,
# each new gap row (2,4,6,8) will hold numeric "activity period" data, and a
# final numeric row 10 is appended after the last text row.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(8).Insert()

# New numeric rows (activity period data)
function Set-RowValues($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

Set-RowValues 2 @(14, 4, 6, 11, 18, 11, 8)
Set-RowValues 4 @(1, 6, 3, 14, 18, 11, 10)
Set-RowValues 6 @(18, 5, 12, 5, 9, 12, 10)
Set-RowValues 8 @(8, 13, 13, 17, 6, 2, 2)
Set-RowValues 10 @(1, 7, 2, 4, 7, 4, 18)

# Row 5 (formerly row 3) referenced "Uberball", which no longer exists as an
# activity. Replace it with "Storm the Castle".
$ws.Range("E5").Value = "Storm the Castle"

# Final cell selection recorded for the sheet view.
$ws.Range("K3").Select()

Write-Output "done"
